# Rename the worksheet "Property1" -> "DataNode" so the sheet's identity
# lines up with the unified DataNode/DataTable/Entity naming scheme used
# across the config resources.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DataNode"

# Reflect where the editor was last positioned on the sheet after the
# rename (matches the saved selection in the authored workbook).
$ws.Activate()
$ws.Range("C24").Select()
